$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.804.18'
$ws.Range("E2").Value = '  +7.33%  '

$ws.Range("D3").Value = '1.744.91'
$ws.Range("E3").Value = '  +3.94%  '

$ws.Range("D4").Value = "'1.001"

$ws.Range("D5").Value = "'335.43"
$ws.Range("E5").Value = '  +2.02%  '

$ws.Range("D6").Value = "'0.9990"
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("D7").Value = "'0.3746"
$ws.Range("E7").Value = '  +2.29%  '

$ws.Range("D8").Value = "'48.45"
$ws.Range("E8").Value = '  +3.40%  '

$ws.Range("D9").Value = "'0.3389"
$ws.Range("E9").Value = '  +4.04%  '

$ws.Range("E10").Value = '  +3.95%  '

$ws.Range("E11").Value = '  +5.58%  '

$ws.Range("D12").Value = "'0.9992"
$ws.Range("E12").Value = '  +0.02%  '

$ws.Range("D13").Value = "'6.419"
$ws.Range("E13").Value = '  +5.30%  '

$ws.Range("D14").Value = "'20.53"
$ws.Range("E14").Value = '  +4.36%  '

$ws.Range("D15").Value = "'7.089"
$ws.Range("E15").Value = '  +6.58%  '

$ws.Range("D16").Value = '1.746.07'
$ws.Range("E16").Value = '  +4.17%  '

$ws.Range("D17").Value = "'0.00001082"
$ws.Range("E17").Value = '  +2.99%  '

$ws.Range("D18").Value = "'0.06738"
$ws.Range("E18").Value = '  +2.31%  '

$ws.Range("D19").Value = "'82.83"
$ws.Range("E19").Value = '  +4.83%  '

$ws.Range("D20").Value = "'0.9983"

$ws.Range("E21").Value = '  +5.21%  '

$ws.Range("D22").Value = "'6.226"
$ws.Range("E22").Value = '  +4.98%  '

$ws.Range("E23").Value = '  -0.48%  '

$ws.Range("D24").Value = '26.784.22'
$ws.Range("E24").Value = '  +7.32%  '

$ws.Range("D25").Value = "'2.465"
$ws.Range("E25").Value = '  +0.84%  '

$ws.Range("D26").Value = "'1.480"
$ws.Range("E26").Value = '  +24.48%  '

$ws.Range("D27").Value = "'2.438"
$ws.Range("E27").Value = '  +0.65%  '

$ws.Range("D28").Value = "'151.86"
$ws.Range("E28").Value = '  +2.55%  '

$ws.Range("D29").Value = "'19.70"
$ws.Range("E29").Value = '  +4.80%  '

$ws.Range("D30").Value = '1.941.71'
$ws.Range("E30").Value = '  +4.25%  '

$ws.Range("D31").Value = "'132.66"
$ws.Range("E31").Value = '  +5.43%  '

$ws.Range("D32").Value = "'4.121"
$ws.Range("E32").Value = '  +1.05%  '

$ws.Range("D33").Value = "'6.065"
$ws.Range("E33").Value = '  +4.70%  '

$ws.Range("D34").Value = "'0.08654"
$ws.Range("E34").Value = '  +1.92%  '

$ws.Range("D35").Value = "'1.694"
$ws.Range("E35").Value = '  +3.01%  '

$ws.Range("E36").Value = '  +4.94%  '

$ws.Range("D37").Value = "'5.445"
$ws.Range("E37").Value = '  +4.90%  '

$ws.Range("D38").Value = "'0.02357"
$ws.Range("E38").Value = '  +4.40%  '

$ws.Range("D39").Value = "'0.2186"
$ws.Range("E39").Value = '  +4.28%  '

$ws.Range("D40").Value = "'0.06277"
$ws.Range("E40").Value = '  +3.95%  '

$ws.Range("D41").Value = "'8.534"
$ws.Range("E41").Value = '  +3.27%  '

$ws.Range("D42").Value = "'1.224"
$ws.Range("E42").Value = '  -0.74%  '

$ws.Range("D43").Value = "'0.6273"
$ws.Range("E43").Value = '  +5.10%  '

$ws.Range("D44").Value = "'14.33"
$ws.Range("E44").Value = '  +4.57%  '

$ws.Range("D45").Value = "'0.9986"
$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("D46").Value = "'3.931"
$ws.Range("E46").Value = '  +2.27%  '

$ws.Range("D47").Value = "'0.6094"
$ws.Range("E47").Value = '  +6.09%  '

$ws.Range("D48").Value = "'129.47"
$ws.Range("E48").Value = '  +3.04%  '

$ws.Range("E49").Value = '  +5.29%  '

$ws.Range("D50").Value = "'0.07224"
$ws.Range("E50").Value = '  +2.96%  '

$ws.Range("D51").Value = "'77.88"
$ws.Range("E51").Value = '  +4.09%  '
